# Update row 3 ("R") values on both the OFF and DEF sheets
# to reflect Week 15 logged data and Week 16 simulated data.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 154
$wsOff.Range("C3").Value = 112
$wsOff.Range("D3").Value = 33
$wsOff.Range("E3").Value = 14
$wsOff.Range("F3").Value = 3

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 220
$wsDef.Range("C3").Value = 139
$wsDef.Range("D3").Value = 61
$wsDef.Range("E3").Value = 22
$wsDef.Range("F3").Value = 7
